$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update item-wise yesterday sales data (columns D:J) to reflect corrected
# item ordering / UOM / sales figures used for the cumulative sales target update.
$data = @(
    ,@(2, 'Augment 375 Tablet 18''s', '18''s', 91, 269.87, 24558.17, 24558.17, 0)
    ,@(3, 'Augment 625 Tablet 18''s', '18''s', 320, 337.33, 107945.6, 107945.6, 0)
    ,@(4, 'Augment 1.2g IV Injection 1''s', '1''s', 127, 206.15, 26181.05, 26181.05, 0)
    ,@(5, 'Augment 100ml PFS', '1''s', 173, 131.18, 22694.14, 22694.14, 0)
    ,@(6, 'Augment 1gm Tablet 12''s', '12''s', 84, 269.87, 22669.08, 22669.08, 0)
    ,@(7, 'Biltin 20mg Tablet 20''s', '20''s', 1880, 224.89, 422793.2, 422793.2, 0)
    ,@(8, 'Bonflex FC Tablet', '40 ''s', 4, 239.88, 959.52, 959.52, 0)
    ,@(9, 'Desodin 5mg Tablet', '50 ''s', 27, 112.45, 3036.15, 3036.15, 0)
    ,@(10, 'Dinafex 60mg FC Tablet 40''s', '40''s', 86, 104.95, 9025.700000000001, 9025.700000000001, 0)
    ,@(11, 'Dinafex 180mg FC Tablet 40''s', '40''s', 112, 299.85, 33583.2, 33583.2, 0)
    ,@(12, 'Dinafex 50ml Suspension', '50 ml', 1482, 35.98, 53322.36, 53322.36, 0)
    ,@(13, 'Dinafex 120mg FC Tablet 40''s', '40''s', 746, 239.88, 178950.48, 178950.48, 0)
    ,@(14, 'Dorenta 100ml Syrup', '100 ml', 1146, 29.99, 34368.54, 34368.54, 0)
    ,@(15, 'Etorix 90mg Tablet 40''s', '40''s', 2467, 359.82, 887675.9399999999, 887675.9399999999, 0)
    ,@(16, 'Etorix 120mg Tablet 30''s', '30''s', 2481, 314.84, 781118.04, 781118.04, 0)
    ,@(17, 'Etorix 60mg Tablet 50''s', '50''s', 1151, 262.37, 301987.87, 301987.87, 0)
    ,@(18, 'Fenobac 5mg Tablet', '50 ''s', 28, 168.67, 4722.76, 4722.76, 0)
    ,@(19, 'Fenobac 10mg Tablet', '30 ''s', 99, 179.91, 17811.09, 17811.09, 0)
    ,@(20, 'Flucloxin 100ml Dry Suspension', '100 ml', 4804, 45.73, 219686.92, 187035.7, 32651.22)
    ,@(21, 'Flucloxin 500mg Capsule 40''s', '40''s', 2717, 314.84, 855420.28, 805990.4, 49429.88)
    ,@(22, 'Flucloxin 250mg Capsule', '100 ''s', 244, 432.53, 105537.32, 105537.32, 0)
    ,@(23, 'Flucloxin 500mg IM/IV Injection', '1''s', 134, 33.94, 4547.96, 4547.96, 0)
    ,@(24, 'Ketonic 10mg Tablet - 30''s', '30''s', 213, 224.89, 47901.57, 47901.57, 0)
    ,@(25, 'Ketonic 60mg IM Injection', '1 ''s', 679, 37.61, 25537.19, 25537.19, 0)
    ,@(26, 'Ketonic 30mg IM/IV Injection', '1 ''s', 1574, 22.83, 35934.42, 35934.42, 0)
    ,@(27, 'Levomax 750mg Tablet - 10''s', '10 ''s', 2, 149.93, 299.86, 299.86, 0)
    ,@(28, 'Levomax 500mg Tablet - 20''s', '20 ''s', 188, 224.89, 42279.32, 42279.32, 0)
    ,@(29, 'Lindamax 150mg Capsule', '30''s', 12, 179.91, 2158.92, 2158.92, 0)
    ,@(30, 'Lindamax Plus 10gm Gel', '10gm', 95, 89.95999999999999, 8546.200000000001, 8546.200000000001, 0)
    ,@(31, 'Lindamax 25ml Lotion', '1''s', 70, 93.7, 6559, 6559, 0)
    ,@(32, 'Lindamax 300mg Capsule', '30''s', 49, 337.33, 16529.17, 16529.17, 0)
    ,@(33, 'Mebidal Tablet', '200 ''s', 20, 449.78, 8995.6, 8995.6, 0)
    ,@(34, 'Nabumet 750mg FC Tab 24''s', '24''s', 121, 395.8, 47891.8, 47891.8, 0)
    ,@(35, 'Nabumet 500mg FC Tab 30''s', '30''s', 971, 337.33, 327547.43, 327547.43, 0)
    ,@(36, 'Naprox Plus 375mg Tablet - 30''s', '30 ''s', 242, 247.38, 59865.96, 59865.96, 0)
    ,@(37, 'Naprox 15gm Gel', '15 gm', 56, 52.47, 2938.32, 2938.32, 0)
    ,@(38, 'Naprox 500mg Tablet', '50 ''s', 558, 412.3, 230063.4, 229651.1, 412.3)
    ,@(39, 'Naprox Plus 500mg Tablet - 36''s', '36''s', 630, 350.82, 221016.6, 221016.6, 0)
    ,@(40, 'Naprox 50ml Suspension', '50 ml', 308, 52.47, 16160.76, 16160.76, 0)
    ,@(41, 'Naprox 250mg Tablet', ' 50 ''s', 222, 262.37, 58246.14, 58246.14, 0)
    ,@(42, 'Ontin 10mg Tablet', '100 ''s', 32, 224.89, 7196.48, 7196.48, 0)
    ,@(43, 'Ontin 60ml Syrup', '60 ml', 346, 22.49, 7781.54, 7781.54, 0)
    ,@(44, 'Oradin 60ml Suspension', '60 ml', 1446, 22.49, 32520.54, 32340.62, 179.92)
    ,@(45, 'Oradin FT 10mg Tablet', '40 ''s', 41, 119.94, 4917.54, 4917.54, 0)
    ,@(46, 'Oradin 10mg Tablet', '100 ''s', 1223, 299.85, 366716.55, 366716.55, 0)
    ,@(47, 'Osticare FC Tab Container 30''s', '30''s', 12, 269.87, 3238.44, 3238.44, 0)
    ,@(48, 'Osticare Tablet 30''s', '30''s', 287, 269.87, 77452.69, 77452.69, 0)
    ,@(49, 'Paino 100mg Tablet', '100 ''s', 84, 299.85, 25187.4, 25187.4, 0)
    ,@(50, 'Quinox DS 60ml Pellets for Suspension', '60 ml', 1351, 74.95999999999999, 101270.96, 90701.60000000001, 10569.36)
    ,@(51, 'Quinox 750mg Tablet', '10 ''s', 8, 134.93, 1079.44, 1079.44, 0)
    ,@(52, 'Quinox 250mg Tablet', '30 ''s', 18, 191.15, 3440.7, 3440.7, 0)
    ,@(53, 'Quinox 500mg Tablet - 20''s', '20 ''s', 87, 224.89, 19565.43, 13268.51, 6296.92)
    ,@(54, 'Quinox 500mg Tablet (40''s)', '40 ''s', 251, 449.78, 112894.78, 107497.42, 5397.36)
    ,@(55, 'Quinox 100ml IV Infusion', '1''s', 10, 74.95999999999999, 749.6, 749.6, 0)
    ,@(56, 'Rupaday 10mg Tablet 30''s', '30''s', 116, 224.89, 26087.24, 26087.24, 0)
    ,@(57, 'Sk-Mox 500mg Capsule 50''s', '50''s', 957, 228.26, 218444.82, 218444.82, 0)
    ,@(58, 'Sk-Mox DS 100ml Dry Suspension', '100 ml', 50, 51.12, 2556, 2556, 0)
    ,@(59, 'Sk-Mox 250mg Capsule', '100 ''s', 37, 266.12, 9846.440000000001, 9846.440000000001, 0)
    ,@(60, 'Sk-Mox 100ml Dry Suspension', '100 ml', 739, 35.68, 26367.52, 26367.52, 0)
    ,@(61, 'Sk-Mox 15ml P/D', '15 ml', 173, 22.64, 3916.72, 3916.72, 0)
    ,@(62, 'Stiba 10mg Tablet - 30''s', '30''s', 219, 179.91, 39400.29, 39400.29, 0)
    ,@(63, 'Stiba 30ml Syrup', '30ml', 30, 44.98, 1349.4, 1349.4, 0)
    ,@(64, 'Sulidac 200 Tablet 20''s', '20''s', 266, 142.43, 37886.38, 37886.38, 0)
    ,@(65, 'Sulidac 100mg Tablet 50''s', '50''s', 82, 187.41, 15367.62, 15367.62, 0)
    ,@(66, 'Tenoxim 20mg Tablet', '30''s', 76, 179.91, 13673.16, 13673.16, 0)
    ,@(67, 'Timothy 5mg IM/IV', '5''s', 96, 56.22, 5397.12, 5397.12, 0)
    ,@(68, 'Timothy 50mg Tablet', '50 ''s', 267, 149.93, 40031.31, 40031.31, 0)
    ,@(69, 'Tojak 5mg Tablet 10''s', '10''S', 294, 359.82, 105787.08, 105787.08, 0)
    ,@(70, 'Toperin 50mg Tablet - 60''s', '60 ''s', 42, 134.93, 5667.06, 5667.06, 0)
    ,@(71, 'Toti 100ml Syrup', '100 ml', 5336, 41.23, 220003.28, 201326.09, 18677.19)
    ,@(72, 'Toti Tablet', '100 ''s', 174, 187.41, 32609.34, 32609.34, 0)
    ,@(73, 'Visomox 400mg FC Tablet 10''s', '10''s', 114, 299.85, 34182.9, 34182.9, 0)
    ,@(74, 'Volmax SR 100mg Capsule', '60 ''s', 10, 134.93, 1349.3, 1349.3, 0)
    ,@(75, 'Xenthol 30 Cream', '15gm', 82, 29.99, 2459.18, 2459.18, 0)
    ,@(76, 'Zithrox 500mg Tablet - 12''s', '12 ''s', 1562, 314.84, 491780.08, 465963.2, 25816.88)
    ,@(77, 'Zithrox 250mg Tablet - 12''s', '12''s', 123, 224.89, 27661.47, 27661.47, 0)
    ,@(78, 'Zithrox 20ml Powder for Suspension', '20ml', 922, 71.20999999999999, 65655.62, 62949.64, 2705.98)
    ,@(79, 'Zithrox 50ml Powder for Suspension', '50ml', 188, 143.93, 27058.84, 27058.84, 0)
    ,@(80, 'Zithrox 35ml Dry Suspension', '35ml', 1616, 104.95, 169599.2, 162672.5, 6926.7)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value2 = $item[1]
    $ws.Cells.Item($r, 5).Value2 = $item[2]
    $ws.Cells.Item($r, 6).Value2 = $item[3]
    $ws.Cells.Item($r, 7).Value2 = $item[4]
    $ws.Cells.Item($r, 8).Value2 = $item[5]
    $ws.Cells.Item($r, 9).Value2 = $item[6]
    $ws.Cells.Item($r, 10).Value2 = $item[7]
}
